# Apply the cryptos.xlsx price/volume refresh described by the commit diff.
# Values are written cell-by-cell to exactly match the target content.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new text would otherwise be auto-parsed by Excel as a number
# (e.g. "1.001") are first forced to Text format so they are stored as
# inline/shared strings, matching the original t="inlineStr" cells.

$ws.Range("D2").Value = "30.543.32"
$ws.Range("E2").Value = "  +0.41%  "
$ws.Range("D3").Value = "1.916.45"
$ws.Range("E3").Value = "  -0.06%  "
$ws.Range("E4").Value = "  +0.13%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "244.09"
$ws.Range("E5").Value = "  +1.32%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.001"
$ws.Range("E6").Value = "  +0.08%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4905"
$ws.Range("E7").Value = "  +4.46%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2900"
$ws.Range("E8").Value = "  +1.36%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06703"
$ws.Range("E9").Value = "  -3.83%  "
$ws.Range("B10").Value = "Litecoin"
$ws.Range("C10").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "105.82"
$ws.Range("E10").Value = "  -0.83%  "
$ws.Range("B11").Value = "Solana"
$ws.Range("C11").Value = "https://coinranking.com/coin/zNZHO_Sjf+solana-sol"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "18.66"
$ws.Range("E11").Value = "  +2.02%  "
$ws.Range("D12").Value = "1.918.32"
$ws.Range("E12").Value = "  +0.06%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.07632"
$ws.Range("E13").Value = "  -0.27%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.263"
$ws.Range("E14").Value = "  +1.35%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.6637"
$ws.Range("E15").Value = "  +0.75%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "271.55"
$ws.Range("E16").Value = "  -6.94%  "
$ws.Range("D17").Value = "30.526.16"
$ws.Range("E17").Value = "  +0.35%  "
$ws.Range("E18").Value = "  +0.10%  "
$ws.Range("B19").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C19").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D19").Value = "2.170.43"
$ws.Range("E19").Value = "  -0.01%  "
$ws.Range("B20").Value = "ShibaInu"
$ws.Range("C20").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.000007513"
$ws.Range("E20").Value = "  -1.41%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "12.78"
$ws.Range("E21").Value = "  -1.81%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.489"
$ws.Range("E22").Value = "  +4.75%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "1.002"
$ws.Range("E23").Value = "  +0.00%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "6.397"
$ws.Range("E24").Value = "  +3.31%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "9.376"
$ws.Range("E25").Value = "  +1.11%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "164.07"
$ws.Range("E26").Value = "  -2.54%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "19.91"
$ws.Range("E27").Value = "  -7.43%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.101"
$ws.Range("E28").Value = "  +2.59%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.1050"
$ws.Range("E29").Value = "  -2.66%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.407"
$ws.Range("E30").Value = "  +2.72%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.111"
$ws.Range("E31").Value = "  -0.80%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.033"
$ws.Range("E32").Value = "  +1.46%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.04969"
$ws.Range("E33").Value = "  -2.07%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.7229"
$ws.Range("E34").Value = "  -2.54%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.131"
$ws.Range("E35").Value = "  -1.18%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.001"
$ws.Range("E36").Value = "  +0.09%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.729"
$ws.Range("E37").Value = "  -0.29%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.02028"
$ws.Range("E38").Value = "  +0.28%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.671"
$ws.Range("E39").Value = "  -0.39%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "110.38"
$ws.Range("E40").Value = "  +1.65%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.008"
$ws.Range("E41").Value = "  -2.66%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.4416"
$ws.Range("E42").Value = "  +4.67%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.8641"
$ws.Range("E43").Value = "  -1.29%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "5.861"
$ws.Range("E44").Value = "  +0.11%  "
$ws.Range("E45").Value = "  +0.09%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "67.63"
$ws.Range("E46").Value = "  -0.09%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "7.219"
$ws.Range("E47").Value = "  +0.45%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "9.261"
$ws.Range("E48").Value = "  +0.16%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.1244"
$ws.Range("E49").Value = "  +2.86%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "47.30"
$ws.Range("E50").Value = "  -12.09%  "
$ws.Range("B51").Value = "WOONetwork"
$ws.Range("C51").Value = "https://coinranking.com/coin/k-J3YwacF+woonetwork-woo"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.2492"
$ws.Range("E51").Value = "  +4.38%  "
